# Fix units of 'Est Read Pairs' to be in millions like other sheets.
# Only for round 1 (rows 2-13 of Sheet1, the "Slide #" groups for Round 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Rename the header cell AA1 from "Est Read Pairs (million) " to "Est Read Pairs"
$ws.Range("AA1").Value = "Est Read Pairs"

# 2. Update the "Est Read Pairs" formulas for round 1 rows (2-13) so the
#    result is expressed in actual read pairs (x1000 more) instead of in
#    thousands, matching the "(million)" -> "" unit relabeling.
#    AA2 holds its own (non-shared) formula.
$ws.Range("AA2").Formula = "=Z2*5000*50000"
#    AA3:AA13 share one formula (anchored at AA3); setting the whole range
#    at once keeps them as a single shared-formula group.
$ws.Range("AA3:AA13").Formula = "=Z3*5000*50000"

# 3. Update the active selection to AA1 (matches the saved view state).
$null = $ws.Range("AA1").Select()
